$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new date column "28-ago" as column BD with the family totals for that date.
$ws.Range("BD1").Value = "28-ago"

$ws.Range("BD2").Value = 18
$ws.Range("BD3").Value = 10
$ws.Range("BD4").Value = 10
$ws.Range("BD5").Value = 13
$ws.Range("BD6").Value = 14
$ws.Range("BD7").Value = 18
$ws.Range("BD8").Value = 10
$ws.Range("BD9").Value = 20
$ws.Range("BD10").Value = 32
$ws.Range("BD11").Value = 19

# Match the formatting of the adjacent columns (header text style, data number style)
$ws.Range("BD1").NumberFormat = $ws.Range("BC1").NumberFormat
$ws.Range("BD2:BD11").NumberFormat = $ws.Range("BC2:BC11").NumberFormat
$ws.Range("BD2:BD11").HorizontalAlignment = $ws.Range("BC2:BC11").HorizontalAlignment

$ws.Range("BB16").Select()
